# chore: merged development into webtool
#
# The upstream "development" branch re-ran the PLZ -> nearest-weather-station
# geocoding match against a fresh batch of cities (Dresden, Koeln, Aachen,
# Kassel, Berlin, Muenchen, Hamburg, Freiburg, Braunschweig) instead of the
# old Dresden-only PLZ list, and the merge drops the two now-duplicated
# trailer rows (12 & 13) that had been left over from a previous paste.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A (PLZ) -----------------------------------------------------
# The new PLZ values were pasted in as numbers (Excel autoconverted the
# first one, which happened to look like a date, into a date serial), so
# the cells keep their numeric literal but get reformatted as Text
# afterwards (format applied retroactively does not convert the stored
# type). Only row 7 ("01067") is typed back in with its leading zero, so
# it stays a genuine text value.
$ws.Range("A2").Value = 44137
$ws.Range("A2").NumberFormat = "@"

$ws.Range("A3").Value = 50667
$ws.Range("A3").NumberFormat = "@"

$ws.Range("A4").Value = 52064
$ws.Range("A4").NumberFormat = "@"

$ws.Range("A5").Value = 34117
$ws.Range("A5").NumberFormat = "@"

$ws.Range("A6").Value = 10115
$ws.Range("A6").NumberFormat = "@"

$ws.Range("A7").NumberFormat = "@"
$ws.Range("A7").Value = "01067"

$ws.Range("A8").Value = 80331
$ws.Range("A8").NumberFormat = "@"

$ws.Range("A9").Value = 20354
$ws.Range("A9").NumberFormat = "@"

$ws.Range("A10").Value = 79100
$ws.Range("A10").NumberFormat = "@"

$ws.Range("A11").Value = 38100
$ws.Range("A11").NumberFormat = "@"

# --- Column D (Filename of matched weather-station file) ----------------
$ws.Range("D9").Value = "TRY2015_535485100234_Jahr.dat"
$ws.Range("D6").Value = "TRY2015_525153133939_Jahr.dat"
$ws.Range("D7").Value = "TRY2015_510342136998_Jahr.dat"
$ws.Range("D8").Value = "TRY2015_480091078440_Jahr.dat"
$ws.Range("D10").Value = "TRY2015_481593115227_Jahr.dat"
$ws.Range("D5").Value = "TRY2015_513148094876_Jahr.dat"
$ws.Range("D3").Value = "TRY2015_509319069572_Jahr.dat"
$ws.Range("D2").Value = "TRY2015_515220074856_Jahr.dat"
$ws.Range("D11").Value = "TRY2015_522733105384_Jahr.dat"
$ws.Range("D4").Value = "TRY2015_507755060854_Jahr.dat"

# --- Drop the two leftover duplicate rows from the old list --------------
$ws.Rows("12:13").Delete()

# --- Cosmetic: widen the Filename column & set print setup ---------------
$ws.Columns("D").ColumnWidth = 38.1
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# --- Leave the selection where the author last clicked --------------------
$ws.Range("D5").Select() | Out-Null
